# Auto-generated edit script applying the Midgardsormr_Profits.xlsx diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2629.5
$ws.Range("I17").Value = 2195
$ws.Range("J17").Value = 2669
$ws.Range("K17").Value = 6585
$ws.Range("L17").Value = 8007
$ws.Range("M17").Value = -6417
$ws.Range("N17").Value = -8343
$ws.Range("H18").Value = 1005.5
$ws.Range("I18").Value = 1005.5
$ws.Range("K18").Value = 1005.5
$ws.Range("M18").Value = -721.5
$ws.Range("H38").Value = 1862.5714
$ws.Range("I38").Value = 134.5
$ws.Range("J38").Value = 4166.6665
$ws.Range("K38").Value = 403.5
$ws.Range("L38").Value = 12499.9995
$ws.Range("M38").Value = -31.5
$ws.Range("N38").Value = -13243.9995
$ws.Range("H55").Value = 188.625
$ws.Range("I55").Value = 185.8
$ws.Range("J55").Value = 193.33333
$ws.Range("K55").Value = 185.8
$ws.Range("L55").Value = 193.33333
$ws.Range("M55").Value = 28.19999999999999
$ws.Range("N55").Value = -621.3333299999999
$ws.Range("H62").Value = 6558.0654
$ws.Range("I62").Value = 6048.359
$ws.Range("K62").Value = 6048.359
$ws.Range("M62").Value = -5424.359
$ws.Range("H64").Value = 14449.096
$ws.Range("I64").Value = 18245.715
$ws.Range("J64").Value = 6855.857
$ws.Range("K64").Value = 18245.715
$ws.Range("L64").Value = 6855.857
$ws.Range("M64").Value = -17997.715
$ws.Range("N64").Value = -7351.857
$ws.Range("H65").Value = 6558.0654
$ws.Range("I65").Value = 6048.359
$ws.Range("K65").Value = 30241.795
$ws.Range("M65").Value = -27121.795
$ws.Range("H67").Value = 14449.096
$ws.Range("I67").Value = 18245.715
$ws.Range("J67").Value = 6855.857
$ws.Range("K67").Value = 18245.715
$ws.Range("L67").Value = 6855.857
$ws.Range("M67").Value = -17387.715
$ws.Range("N67").Value = -8571.857
$ws.Range("H70").Value = 6779.364
$ws.Range("J70").Value = 4721.2383
$ws.Range("L70").Value = 14163.7149
$ws.Range("N70").Value = -14703.7149
$ws.Range("H73").Value = 6779.364
$ws.Range("J73").Value = 4721.2383
$ws.Range("L73").Value = 14163.7149
$ws.Range("N73").Value = -16035.7149
$ws.Range("H88").Value = 11818.8
$ws.Range("I88").Value = 1199
$ws.Range("J88").Value = 14473.75
$ws.Range("K88").Value = 1199
$ws.Range("L88").Value = 14473.75
$ws.Range("M88").Value = -793
$ws.Range("N88").Value = -15285.75
$ws.Range("H91").Value = 11818.8
$ws.Range("I91").Value = 1199
$ws.Range("J91").Value = 14473.75
$ws.Range("K91").Value = 1199
$ws.Range("L91").Value = 14473.75
$ws.Range("M91").Value = 205
$ws.Range("N91").Value = -17281.75
$ws.Range("H98").Value = 2391.2046
$ws.Range("J98").Value = 4266
$ws.Range("L98").Value = 4266
$ws.Range("N98").Value = -7262
$ws.Range("H100").Value = 38764.707
$ws.Range("I100").Value = 38764.707
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 38764.707
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -38223.707
$ws.Range("N100").ClearContents()
$ws.Range("H101").Value = 391.66666
$ws.Range("I101").Value = 136.66667
$ws.Range("J101").Value = 646.6667
$ws.Range("K101").Value = 410.00001
$ws.Range("L101").Value = 1940.0001
$ws.Range("M101").Value = 1211.99999
$ws.Range("N101").Value = -5184.0001
$ws.Range("H113").Value = 7248.75
$ws.Range("H122").Value = 2391.2046
$ws.Range("J122").Value = 4266
$ws.Range("L122").Value = 12798
$ws.Range("N122").Value = -17698
$ws.Range("H126").Value = 116666.336
$ws.Range("J126").Value = 130499.5
$ws.Range("L126").Value = 130499.5
$ws.Range("N126").Value = -140379.5
$ws.Range("H132").Value = 2565683.8
$ws.Range("I132").Value = 2565683.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7697051.399999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7694521.399999999
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 107334.664
$ws.Range("I137").Value = 155002
$ws.Range("J137").Value = 12000
$ws.Range("K137").Value = 465006
$ws.Range("L137").Value = 36000
$ws.Range("M137").Value = -462456
$ws.Range("N137").Value = -41100
$ws.Range("H138").Value = 1593.5428
$ws.Range("I138").Value = 1171.4036
$ws.Range("J138").Value = 3444.4614
$ws.Range("K138").Value = 3514.2108
$ws.Range("L138").Value = 10333.3842
$ws.Range("M138").Value = 1625.7892
$ws.Range("N138").Value = -20613.3842
$ws.Range("H141").Value = 1925.909
$ws.Range("I141").Value = 1618.8
$ws.Range("K141").Value = 4856.4
$ws.Range("M141").Value = 323.6000000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 41500
$ws.Range("J24").Value = 41500
$ws.Range("L24").Value = 41500
$ws.Range("N24").Value = -42248
$ws.Range("H32").Value = 17534.754
$ws.Range("I32").Value = 18864.258
$ws.Range("J32").Value = 4999.4287
$ws.Range("K32").Value = 18864.258
$ws.Range("L32").Value = 4999.4287
$ws.Range("M32").Value = -18577.258
$ws.Range("N32").Value = -5573.4287
$ws.Range("H45").Value = 6528.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 6528.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6528.5
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -7282.5
$ws.Range("H63").Value = 1900
$ws.Range("I63").Value = 1900
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1900
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1214
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1900
$ws.Range("I66").Value = 1900
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9500
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6068
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 142867.42
$ws.Range("I74").Value = 154828.2
$ws.Range("K74").Value = 154828.2
$ws.Range("M74").Value = -153954.2
$ws.Range("H77").Value = 142867.42
$ws.Range("I77").Value = 154828.2
$ws.Range("K77").Value = 774141
$ws.Range("M77").Value = -769773
$ws.Range("H88").Value = 4862.846
$ws.Range("I88").Value = 321.7
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 321.7
$ws.Range("L88").Value = 20000
$ws.Range("M88").Value = 84.30000000000001
$ws.Range("N88").Value = -20812
$ws.Range("H91").Value = 4862.846
$ws.Range("I91").Value = 321.7
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 321.7
$ws.Range("L91").Value = 20000
$ws.Range("M91").Value = 1082.3
$ws.Range("N91").Value = -22808
$ws.Range("H100").Value = 41500
$ws.Range("J100").Value = 41500
$ws.Range("L100").Value = 41500
$ws.Range("N100").Value = -43664
$ws.Range("H122").Value = 1623.283
$ws.Range("I122").Value = 1612.9592
$ws.Range("K122").Value = 4838.8776
$ws.Range("M122").Value = -2388.8776
$ws.Range("H132").Value = 1504.1233
$ws.Range("I132").Value = 1218.2131
$ws.Range("J132").Value = 2957.5
$ws.Range("K132").Value = 3654.6393
$ws.Range("L132").Value = 8872.5
$ws.Range("M132").Value = -1124.6393
$ws.Range("N132").Value = -13932.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1553.25
$ws.Range("I86").Value = 1490.3077
$ws.Range("J86").Value = 1826
$ws.Range("K86").Value = 1490.3077
$ws.Range("L86").Value = 1826
$ws.Range("M86").Value = -367.3077000000001
$ws.Range("N86").Value = -4072
$ws.Range("H89").Value = 1553.25
$ws.Range("I89").Value = 1490.3077
$ws.Range("J89").Value = 1826
$ws.Range("K89").Value = 7451.538500000001
$ws.Range("L89").Value = 9130
$ws.Range("M89").Value = -1835.538500000001
$ws.Range("N89").Value = -20362
$ws.Range("H94").Value = 2816.5
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 3624.75
$ws.Range("K94").Value = 1200
$ws.Range("L94").Value = 3624.75
$ws.Range("M94").Value = -749
$ws.Range("N94").Value = -4526.75
$ws.Range("H99").Value = 2553.3333
$ws.Range("I99").Value = 2211.4285
$ws.Range("K99").Value = 2211.4285
$ws.Range("M99").Value = -713.4285
$ws.Range("H105").Value = 4108.933
$ws.Range("I105").Value = 4481.4346
$ws.Range("K105").Value = 4481.4346
$ws.Range("M105").Value = -2734.4346
$ws.Range("H107").Value = 23364.5
$ws.Range("I107").Value = 37442.645
$ws.Range("K107").Value = 37442.645
$ws.Range("M107").Value = -35522.645
$ws.Range("H133").Value = 99780
$ws.Range("J133").Value = 99780
$ws.Range("L133").Value = 99780
$ws.Range("N133").Value = -109900
$ws.Range("H134").Value = 1501.5333
$ws.Range("I134").Value = 1308.9857
$ws.Range("K134").Value = 3926.9571
$ws.Range("M134").Value = -1391.9571
$ws.Range("H137").Value = 121670
$ws.Range("J137").Value = 121670
$ws.Range("L137").Value = 121670
$ws.Range("N137").Value = -131870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1142.2307
$ws.Range("I16").Value = 1164.3334
$ws.Range("K16").Value = 1164.3334
$ws.Range("M16").Value = -877.3334
$ws.Range("H31").Value = 9094164
$ws.Range("I31").Value = 20001996
$ws.Range("J31").Value = 4304.1665
$ws.Range("K31").Value = 20001996
$ws.Range("L31").Value = 4304.1665
$ws.Range("M31").Value = -20001701
$ws.Range("N31").Value = -4894.1665
$ws.Range("H34").Value = 9094164
$ws.Range("I34").Value = 20001996
$ws.Range("J34").Value = 4304.1665
$ws.Range("K34").Value = 20001996
$ws.Range("L34").Value = 4304.1665
$ws.Range("M34").Value = -20001794
$ws.Range("N34").Value = -4708.1665
$ws.Range("H58").Value = 979.5862
$ws.Range("I58").Value = 1022.875
$ws.Range("K58").Value = 1022.875
$ws.Range("M58").Value = -819.875
$ws.Range("H105").Value = 2188.1667
$ws.Range("I105").Value = 1173.3334
$ws.Range("J105").Value = 3203
$ws.Range("K105").Value = 1173.3334
$ws.Range("L105").Value = 3203
$ws.Range("M105").Value = 573.6666
$ws.Range("N105").Value = -6697
$ws.Range("H107").Value = 4649.1665
$ws.Range("J107").Value = 4649.1665
$ws.Range("L107").Value = 4649.1665
$ws.Range("N107").Value = -8489.166499999999
$ws.Range("H113").Value = 1142.2307
$ws.Range("I113").Value = 1164.3334
$ws.Range("K113").Value = 1164.3334
$ws.Range("M113").Value = 1005.6666
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H131").Value = 105665.664
$ws.Range("J131").Value = 105665.664
$ws.Range("L131").Value = 105665.664
$ws.Range("N131").Value = -115745.664
$ws.Range("H132").Value = 20432.656
$ws.Range("I132").Value = 27441.738
$ws.Range("K132").Value = 82325.21400000001
$ws.Range("M132").Value = -79795.21400000001
$ws.Range("H134").Value = 1933.5172
$ws.Range("J134").Value = 3495.1667
$ws.Range("L134").Value = 10485.5001
$ws.Range("N134").Value = -15555.5001
$ws.Range("H136").Value = 979.5862
$ws.Range("I136").Value = 1022.875
$ws.Range("K136").Value = 3068.625
$ws.Range("M136").Value = -518.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 271.1111
$ws.Range("I8").Value = 271.1111
$ws.Range("K8").Value = 813.3333
$ws.Range("M8").Value = -674.3333
$ws.Range("H69").Value = 4950
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -11189
$ws.Range("N69").Value = -16622
$ws.Range("H72").Value = 4950
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -31944
$ws.Range("N72").Value = -53112
$ws.Range("H80").Value = 4875.615
$ws.Range("J80").Value = 4870.72
$ws.Range("L80").Value = 14612.16
$ws.Range("N80").Value = -16484.16
$ws.Range("H83").Value = 4875.615
$ws.Range("J83").Value = 4870.72
$ws.Range("L83").Value = 43836.48
$ws.Range("N83").Value = -53196.48
$ws.Range("H97").Value = 2209.5
$ws.Range("I97").Value = 279.5
$ws.Range("J97").Value = 4139.5
$ws.Range("K97").Value = 838.5
$ws.Range("L97").Value = 12418.5
$ws.Range("M97").Value = -342.5
$ws.Range("N97").Value = -13410.5
$ws.Range("H103").Value = 655.5
$ws.Range("J103").Value = 1215
$ws.Range("L103").Value = 3645
$ws.Range("N103").Value = -5403
$ws.Range("H104").Value = 2800
$ws.Range("H105").Value = 7374.375
$ws.Range("I105").Value = 7374.375
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 22123.125
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -19502.125
$ws.Range("N105").ClearContents()
$ws.Range("H109").Value = 1690.9166
$ws.Range("J109").Value = 2856.8333
$ws.Range("L109").Value = 8570.499899999999
$ws.Range("N109").Value = -10650.4999
$ws.Range("H111").Value = 10296
$ws.Range("I111").Value = 11480
$ws.Range("K111").Value = 34440
$ws.Range("M111").Value = -31373
$ws.Range("H113").Value = 2051.5715
$ws.Range("I113").Value = 594.25
$ws.Range("K113").Value = 1782.75
$ws.Range("M113").Value = 387.25
$ws.Range("H114").Value = 1036
$ws.Range("J114").Value = 1051.1
$ws.Range("L114").Value = 3153.3
$ws.Range("N114").Value = -9661.299999999999
$ws.Range("H117").Value = 1755
$ws.Range("I117").Value = 2193.1667
$ws.Range("J117").Value = 878.6667
$ws.Range("K117").Value = 6579.500100000001
$ws.Range("L117").Value = 2636.0001
$ws.Range("M117").Value = -3137.500100000001
$ws.Range("N117").Value = -9520.000100000001
$ws.Range("H139").Value = 3791.8333
$ws.Range("I139").Value = 3994
$ws.Range("J139").Value = 2376.6667
$ws.Range("K139").Value = 11982
$ws.Range("L139").Value = 7130.000100000001
$ws.Range("M139").Value = -6842
$ws.Range("N139").Value = -17410.0001
$ws.Range("H140").Value = 2995.6316
$ws.Range("I140").Value = 2995.6316
$ws.Range("K140").Value = 8986.8948
$ws.Range("M140").Value = -3806.8948
$ws.Range("H141").Value = 4254.2
$ws.Range("I141").Value = 4200.9287
$ws.Range("K141").Value = 12602.7861
$ws.Range("M141").Value = -7422.786100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 491.2
$ws.Range("I107").Value = 742.375
$ws.Range("K107").Value = 742.375
$ws.Range("M107").Value = 1177.625
$ws.Range("H113").Value = 2279.2
$ws.Range("I113").Value = 2032.6666
$ws.Range("J113").Value = 2649
$ws.Range("K113").Value = 2032.6666
$ws.Range("L113").Value = 2649
$ws.Range("M113").Value = 137.3334
$ws.Range("N113").Value = -6989
$ws.Range("H126").Value = 6229.9
$ws.Range("I126").Value = 2459.8
$ws.Range("K126").Value = 7379.400000000001
$ws.Range("M126").Value = -4909.400000000001
$ws.Range("H132").Value = 1967.8
$ws.Range("I132").Value = 1790.3823
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 5371.1469
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -2841.1469
$ws.Range("N132").Value = -29060
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2334.1035
$ws.Range("I7").Value = 2203.1785
$ws.Range("K7").Value = 2203.1785
$ws.Range("M7").Value = -2091.1785
$ws.Range("H40").Value = 2100.8
$ws.Range("I40").Value = 2100.8
$ws.Range("K40").Value = 2100.8
$ws.Range("M40").Value = -1964.8
$ws.Range("H46").Value = 5471.45
$ws.Range("J46").Value = 6884.2144
$ws.Range("L46").Value = 6884.2144
$ws.Range("N46").Value = -7260.2144
$ws.Range("H68").Value = 2906.5151
$ws.Range("I68").Value = 2531.1538
$ws.Range("J68").Value = 3150.5
$ws.Range("K68").Value = 2531.1538
$ws.Range("L68").Value = 3150.5
$ws.Range("M68").Value = -1782.1538
$ws.Range("N68").Value = -4648.5
$ws.Range("H71").Value = 2906.5151
$ws.Range("I71").Value = 2531.1538
$ws.Range("J71").Value = 3150.5
$ws.Range("K71").Value = 12655.769
$ws.Range("L71").Value = 15752.5
$ws.Range("M71").Value = -8911.769
$ws.Range("N71").Value = -23240.5
$ws.Range("H93").Value = 2027.7941
$ws.Range("I93").Value = 2161.8635
$ws.Range("J93").Value = 1782
$ws.Range("K93").Value = 2161.8635
$ws.Range("L93").Value = 1782
$ws.Range("M93").Value = -913.8634999999999
$ws.Range("N93").Value = -4278
$ws.Range("H100").Value = 2507
$ws.Range("I100").Value = 2137.25
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2137.25
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1596.25
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 2901.6562
$ws.Range("I122").Value = 2884.4482
$ws.Range("K122").Value = 8653.3446
$ws.Range("M122").Value = -6203.3446
$ws.Range("H126").Value = 2334.1035
$ws.Range("I126").Value = 2203.1785
$ws.Range("K126").Value = 6609.5355
$ws.Range("M126").Value = -4139.5355
$ws.Range("H132").Value = 2266.75
$ws.Range("I132").Value = 2188.6858
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6566.057400000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4036.057400000001
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 101201.75
$ws.Range("J134").Value = 101201.75
$ws.Range("L134").Value = 101201.75
$ws.Range("N134").Value = -111341.75
$ws.Range("H136").Value = 4371.8276
$ws.Range("I136").Value = 3784.6086
$ws.Range("J136").Value = 6622.8335
$ws.Range("K136").Value = 11353.8258
$ws.Range("L136").Value = 19868.5005
$ws.Range("M136").Value = -8803.825800000001
$ws.Range("N136").Value = -24968.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 43975
$ws.Range("J58").Value = 43975
$ws.Range("L58").Value = 43975
$ws.Range("N58").Value = -44591
$ws.Range("H62").Value = 4666.6665
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 4666.6665
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -36240
$ws.Range("H132").Value = 4824257.5
$ws.Range("I132").Value = 6431236
$ws.Range("J132").Value = 3322.5386
$ws.Range("K132").Value = 19293708
$ws.Range("L132").Value = 9967.6158
$ws.Range("M132").Value = -19291178
$ws.Range("N132").Value = -15027.6158
$ws.Range("H136").Value = 12962.302
$ws.Range("I136").Value = 13964.106
$ws.Range("K136").Value = 41892.318
$ws.Range("M136").Value = -39342.318

Write-Output "Applied 496 cell changes across 8 sheets"
